$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NroSiniestro value in F2 (row 2) with the new claim number.
# A leading apostrophe forces Excel to keep it as text (matching the
# original cell's quote-prefixed text style) instead of converting it
# to a number and losing the leading zero.
$ws.Range("F2").Value = "'0420194406719"

# Update the active selection to match the new cursor position
$ws.Range("H6").Select()
